# Auto-generated Excel COM-interop script
# Applies updated market price data (columns H-N) to specific rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled runner diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 7844.5557
$ws.Range("I74").Value = 5866
$ws.Range("J74").Value = 8833.833000000001
$ws.Range("K74").Value = 5866
$ws.Range("L74").Value = 8833.833000000001
$ws.Range("M74").Value = -4930
$ws.Range("N74").Value = -10705.833
# Row 77
$ws.Range("H77").Value = 7844.5557
$ws.Range("I77").Value = 5866
$ws.Range("J77").Value = 8833.833000000001
$ws.Range("K77").Value = 29330
$ws.Range("L77").Value = 44169.165
$ws.Range("M77").Value = -24650
$ws.Range("N77").Value = -53529.165
# Row 80
$ws.Range("H80").Value = 951.913
$ws.Range("I80").Value = 187.23529
$ws.Range("J80").Value = 3118.5
$ws.Range("K80").Value = 561.70587
$ws.Range("L80").Value = 9355.5
$ws.Range("M80").Value = 436.29413
$ws.Range("N80").Value = -11351.5
# Row 82
$ws.Range("H82").Value = 3045.6667
$ws.Range("I82").Value = 3045.6667
$ws.Range("K82").Value = 9137.000100000001
$ws.Range("M82").Value = -8731.000100000001
# Row 83
$ws.Range("H83").Value = 951.913
$ws.Range("I83").Value = 187.23529
$ws.Range("J83").Value = 3118.5
$ws.Range("K83").Value = 1685.11761
$ws.Range("L83").Value = 28066.5
$ws.Range("M83").Value = 3306.88239
$ws.Range("N83").Value = -38050.5
# Row 85
$ws.Range("H85").Value = 3045.6667
$ws.Range("I85").Value = 3045.6667
$ws.Range("K85").Value = 9137.000100000001
$ws.Range("M85").Value = -7733.000100000001
# Row 112
$ws.Range("H112").Value = 2248.4443
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2248.4443
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6745.3329
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -8961.332900000001
# Row 116
$ws.Range("H116").Value = 9849.125
$ws.Range("I116").Value = 9132.166999999999
$ws.Range("K116").Value = 9132.166999999999
$ws.Range("M116").Value = -5690.166999999999
# Row 138
$ws.Range("H138").Value = 3064.2964
$ws.Range("J138").Value = 3966.0688
$ws.Range("L138").Value = 11898.2064
$ws.Range("N138").Value = -22178.2064

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 3993.75
$ws.Range("I74").Value = 1758.65
$ws.Range("J74").Value = 9581.5
$ws.Range("K74").Value = 1758.65
$ws.Range("L74").Value = 9581.5
$ws.Range("M74").Value = -884.6500000000001
$ws.Range("N74").Value = -11329.5
# Row 77
$ws.Range("H77").Value = 3993.75
$ws.Range("I77").Value = 1758.65
$ws.Range("J77").Value = 9581.5
$ws.Range("K77").Value = 8793.25
$ws.Range("L77").Value = 47907.5
$ws.Range("M77").Value = -4425.25
$ws.Range("N77").Value = -56643.5
# Row 102
$ws.Range("H102").Value = 10870938
$ws.Range("I102").Value = 1435.1364
$ws.Range("K102").Value = 1435.1364
$ws.Range("M102").Value = 186.8635999999999
# Row 110
$ws.Range("H110").Value = 7354475.5
$ws.Range("I110").Value = 11364475
$ws.Range("K110").Value = 11364475
$ws.Range("M110").Value = -11362430
# Row 122
$ws.Range("H122").Value = 2377
$ws.Range("I122").Value = 2921.0908
$ws.Range("J122").Value = 1180
$ws.Range("K122").Value = 8763.2724
$ws.Range("L122").Value = 3540
$ws.Range("M122").Value = -6313.2724
$ws.Range("N122").Value = -8440
# Row 132
$ws.Range("H132").Value = 4353237.5
$ws.Range("I132").Value = 5886235.5
$ws.Range("J132").Value = 9742.75
$ws.Range("K132").Value = 17658706.5
$ws.Range("L132").Value = 29228.25
$ws.Range("M132").Value = -17656176.5
$ws.Range("N132").Value = -34288.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2606.7273
$ws.Range("I20").Value = 2396
$ws.Range("J20").Value = 3091.4
$ws.Range("K20").Value = 2396
$ws.Range("L20").Value = 3091.4
$ws.Range("M20").Value = -2149
$ws.Range("N20").Value = -3585.4
# Row 94
$ws.Range("H94").Value = 2037.7097
$ws.Range("I94").Value = 2290.1667
$ws.Range("K94").Value = 2290.1667
$ws.Range("M94").Value = -1839.1667
# Row 105
$ws.Range("H105").Value = 52646868
$ws.Range("I105").Value = 76943976
$ws.Range("K105").Value = 76943976
$ws.Range("M105").Value = -76942229
# Row 107
$ws.Range("H107").Value = 1732.5834
$ws.Range("I107").Value = 1294.4375
$ws.Range("K107").Value = 1294.4375
$ws.Range("M107").Value = 625.5625
# Row 126
$ws.Range("H126").Value = 47250
$ws.Range("J126").Value = 47250
$ws.Range("L126").Value = 47250
$ws.Range("N126").Value = -57130
# Row 134
$ws.Range("H134").Value = 8244.764999999999
$ws.Range("I134").Value = 7144.067
$ws.Range("K134").Value = 21432.201
$ws.Range("M134").Value = -18897.201

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 469.44446
$ws.Range("J22").Value = 517.8
$ws.Range("L22").Value = 517.8
$ws.Range("N22").Value = -1217.8
# Row 52
$ws.Range("H52").Value = 67854.5
$ws.Range("I52").Value = 67854.5
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 67854.5
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -67560.5
$ws.Range("N52").ClearContents()
# Row 62
$ws.Range("H62").Value = 10114.789
$ws.Range("I62").Value = 7981.4287
$ws.Range("J62").Value = 11359.25
$ws.Range("K62").Value = 7981.4287
$ws.Range("L62").Value = 11359.25
$ws.Range("M62").Value = -7357.4287
$ws.Range("N62").Value = -12607.25
# Row 65
$ws.Range("H65").Value = 10114.789
$ws.Range("I65").Value = 7981.4287
$ws.Range("J65").Value = 11359.25
$ws.Range("K65").Value = 39907.14350000001
$ws.Range("L65").Value = 56796.25
$ws.Range("M65").Value = -36787.14350000001
$ws.Range("N65").Value = -63036.25
# Row 138
$ws.Range("H138").Value = 119999.664
$ws.Range("J138").Value = 119999.664
$ws.Range("L138").Value = 119999.664
$ws.Range("N138").Value = -130279.664

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 9460517
$ws.Range("I4").Value = 11914676
$ws.Range("J4").Value = 6337042
$ws.Range("K4").Value = 35744028
$ws.Range("L4").Value = 19011126
$ws.Range("M4").Value = -35743916
$ws.Range("N4").Value = -19011350
# Row 5
$ws.Range("H5").Value = 296.35715
$ws.Range("I5").Value = 273
$ws.Range("K5").Value = 819
$ws.Range("M5").Value = -707
# Row 39
$ws.Range("H39").Value = 2483.3333
$ws.Range("J39").Value = 4875
$ws.Range("L39").Value = 14625
$ws.Range("N39").Value = -15213
# Row 113
$ws.Range("H113").Value = 3588.75
$ws.Range("I113").Value = 2358
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 7074
$ws.Range("L113").Value = 11997
$ws.Range("M113").Value = -4904
$ws.Range("N113").Value = -16337
# Row 135
$ws.Range("H135").Value = 296.35715
$ws.Range("I135").Value = 273
$ws.Range("K135").Value = 2457
$ws.Range("M135").Value = 78

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 11647.714
$ws.Range("I43").Value = 5385
$ws.Range("J43").Value = 19998
$ws.Range("K43").Value = 5385
$ws.Range("L43").Value = 19998
$ws.Range("M43").Value = -5234
$ws.Range("N43").Value = -20300
# Row 70
$ws.Range("H70").Value = 5228.8237
$ws.Range("I70").Value = 4206.7856
$ws.Range("K70").Value = 4206.7856
$ws.Range("M70").Value = -3936.7856
# Row 73
$ws.Range("H73").Value = 5228.8237
$ws.Range("I73").Value = 4206.7856
$ws.Range("K73").Value = 4206.7856
$ws.Range("M73").Value = -3270.7856
# Row 122
$ws.Range("H122").Value = 5367.222
$ws.Range("I122").Value = 1989.4166
$ws.Range("J122").Value = 12122.833
$ws.Range("K122").Value = 5968.2498
$ws.Range("L122").Value = 36368.499
$ws.Range("M122").Value = -3518.2498
$ws.Range("N122").Value = -41268.499
# Row 132
$ws.Range("H132").Value = 6688
$ws.Range("I132").Value = 4489.2856
$ws.Range("J132").Value = 11818.333
$ws.Range("K132").Value = 13467.8568
$ws.Range("L132").Value = 35454.999
$ws.Range("M132").Value = -10937.8568
$ws.Range("N132").Value = -40514.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 7532.92
$ws.Range("I46").Value = 1220
$ws.Range("J46").Value = 9111.15
$ws.Range("K46").Value = 1220
$ws.Range("L46").Value = 9111.15
$ws.Range("M46").Value = -1032
$ws.Range("N46").Value = -9487.15
# Row 68
$ws.Range("H68").Value = 4374.875
$ws.Range("J68").Value = 5199.8
$ws.Range("L68").Value = 5199.8
$ws.Range("N68").Value = -6697.8
# Row 71
$ws.Range("H71").Value = 4374.875
$ws.Range("J71").Value = 5199.8
$ws.Range("L71").Value = 25999
$ws.Range("N71").Value = -33487
# Row 93
$ws.Range("H93").Value = 1839.3
$ws.Range("I93").Value = 1932.6666
$ws.Range("J93").Value = 999
$ws.Range("K93").Value = 1932.6666
$ws.Range("L93").Value = 999
$ws.Range("M93").Value = -684.6666
$ws.Range("N93").Value = -3495
# Row 137
$ws.Range("H137").Value = 86999
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3559.111
$ws.Range("I122").Value = 3660.25
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 10980.75
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -8530.75
$ws.Range("N122").Value = -13150
# Row 136
$ws.Range("H136").Value = 3761417
$ws.Range("I136").Value = 4465595
$ws.Range("J136").Value = 5800
$ws.Range("K136").Value = 13396785
$ws.Range("L136").Value = 17400
$ws.Range("M136").Value = -13394235
$ws.Range("N136").Value = -22500
